# Add "2022-Q4" fund-holdings sheet right after "总计", shifting the
# existing quarter sheets (2022-Q3 .. 2021-Q1) one position to the right,
# and insert the corresponding summary row into "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for 2022-Q4 and bump the
#    running index in column A for every row that shifts down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# reuse the existing "index column" style (s="2") for the new A2 cell
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 19
$total.Range("D2").Value = 6.01

for ($r = 9; $r -ge 3; $r--) {
    $old = [double]($total.Cells.Item($r, 1).Value2)
    $total.Cells.Item($r, 1).Value = $old + 1
}

# ---------------------------------------------------------------------
# 2) New "2022-Q4" sheet with the fund-holdings detail, placed right
#    after "总计" (pushing 2022-Q3..2021-Q1 one slot later).
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$headerCols = @("B", "C", "D", "E", "F", "G", "H")

# header row formatting/style copied from an existing detail sheet's header
$wb.Worksheets.Item("2022-Q3").Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Range("$($headerCols[$i])1").Value = $headers[$i]
}

$rows = @(
    @(0, "000628", "大成高新技术产业股票A", "46.16", "87.82", "5.82", "2.6865", 5),
    @(1, "010178", "大成企业能力驱动混合A", "31.86", "87.25", "2.87", "0.9144", 9),
    @(2, "008269", "大成睿享混合A", "23.67", "63.04", "2.62", "0.6202", 8),
    @(3, "011066", "大成高新技术产业股票C", "8.86", "87.82", "5.82", "0.5157", 5),
    @(4, "013853", "大成匠心卓越三年持有混合A", "4.09", "68.84", "7.26", "0.2969", 2),
    @(5, "090013", "大成竞争优势混合", "8.41", "62.88", "2.88", "0.2422", 7),
    @(6, "008270", "大成睿享混合C", "8.29", "63.04", "2.62", "0.2172", 8),
    @(7, "009223", "宝盈现代服务业混合A", "3.46", "88.59", "4.40", "0.1522", 10),
    @(8, "015564", "大成弘远回报一年持有混合A", "2.54", "27.63", "4.82", "0.1224", 2),
    @(9, "013463", "大成致远优势一年持有期混合A", "2.63", "66.62", "3.00", "0.0789", 7),
    @(10, "011834", "大成投资严选六月持有混合A", "2.99", "64.26", "2.18", "0.0652", 6),
    @(11, "008303", "宝盈龙头优选股票A", "0.65", "88.21", "3.42", "0.0222", 10),
    @(12, "013854", "大成匠心卓越三年持有混合C", "0.30", "68.84", "7.26", "0.0218", 2),
    @(13, "009224", "宝盈现代服务业混合C", "0.41", "88.59", "4.40", "0.0180", 10),
    @(14, "010179", "大成企业能力驱动混合C", "0.61", "87.25", "2.87", "0.0175", 9),
    @(15, "008304", "宝盈龙头优选股票C", "0.18", "88.21", "3.42", "0.0062", 10),
    @(16, "011835", "大成投资严选六月持有混合C", "0.20", "64.26", "2.18", "0.0044", 6),
    @(17, "015565", "大成弘远回报一年持有混合C", "0.09", "27.63", "4.82", "0.0043", 2),
    @(18, "013464", "大成致远优势一年持有期混合C", "0.14", "66.62", "3.00", "0.0042", 7)
)

$textCols = @("B", "C", "D", "E", "F", "G")

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    # A: numeric running index, styled like the "总计" index column
    $total.Range("A2").Copy()
    $q4.Range("A$r").PasteSpecial(-4122)
    $q4.Range("A$r").Value = $row[0]

    # B..G: text-typed values (kept as text even though some look numeric)
    for ($c = 0; $c -lt $textCols.Length; $c++) {
        $addr = "$($textCols[$c])$r"
        $q4.Range($addr).NumberFormat = "@"
        $q4.Range($addr).Value = $row[$c + 1]
        $q4.Range($addr).Style = "Normal"
    }

    # H: plain numeric rank
    $q4.Range("H$r").Value = $row[7]
}

Write-Host "2022-Q4 sheet inserted; 总计 updated"
